$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "73-9=64"
$tbl.Cell(1, 2).Range.Text = "7+5=12"
$tbl.Cell(1, 3).Range.Text = "34+39=73"
$tbl.Cell(1, 4).Range.Text = "72-64=8"
$tbl.Cell(1, 5).Range.Text = "80-16=64"
$tbl.Cell(2, 1).Range.Text = "3+89=92"
$tbl.Cell(2, 2).Range.Text = "35-8=27"
$tbl.Cell(2, 3).Range.Text = "70-32=38"
$tbl.Cell(2, 4).Range.Text = "62-55=7"
$tbl.Cell(2, 5).Range.Text = "61-4=57"
$tbl.Cell(3, 1).Range.Text = "16+35=51"
$tbl.Cell(3, 2).Range.Text = "92-19=73"
$tbl.Cell(3, 3).Range.Text = "40-29=11"
$tbl.Cell(3, 4).Range.Text = "29+25=54"
$tbl.Cell(3, 5).Range.Text = "14+9=23"
$tbl.Cell(4, 1).Range.Text = "9+12=21"
$tbl.Cell(4, 2).Range.Text = "51-2=49"
$tbl.Cell(4, 3).Range.Text = "26+39=65"
$tbl.Cell(4, 4).Range.Text = "72-45=27"
$tbl.Cell(4, 5).Range.Text = "21-5=16"
$tbl.Cell(5, 1).Range.Text = "15+27=42"
$tbl.Cell(5, 2).Range.Text = "72-45=27"
$tbl.Cell(5, 3).Range.Text = "85-67=18"
$tbl.Cell(5, 4).Range.Text = "77-38=39"
$tbl.Cell(5, 5).Range.Text = "71-7=64"
$tbl.Cell(6, 1).Range.Text = "50-31=19"
$tbl.Cell(6, 2).Range.Text = "45+9=54"
$tbl.Cell(6, 3).Range.Text = "19+49=68"
$tbl.Cell(6, 4).Range.Text = "50-29=21"
$tbl.Cell(6, 5).Range.Text = "38+19=57"
$tbl.Cell(7, 1).Range.Text = "36-18=18"
$tbl.Cell(7, 2).Range.Text = "37+14=51"
$tbl.Cell(7, 3).Range.Text = "28+38=66"
$tbl.Cell(7, 4).Range.Text = "92-74=18"
$tbl.Cell(7, 5).Range.Text = "64-29=35"
$tbl.Cell(8, 1).Range.Text = "39+33=72"
$tbl.Cell(8, 2).Range.Text = "8+83=91"
$tbl.Cell(8, 3).Range.Text = "77-28=49"
$tbl.Cell(8, 4).Range.Text = "85+8=93"
$tbl.Cell(8, 5).Range.Text = "20-15=5"
$tbl.Cell(9, 1).Range.Text = "47-38=9"
$tbl.Cell(9, 2).Range.Text = "37+36=73"
$tbl.Cell(9, 3).Range.Text = "94-65=29"
$tbl.Cell(9, 4).Range.Text = "29+62=91"
$tbl.Cell(9, 5).Range.Text = "24+18=42"
$tbl.Cell(10, 1).Range.Text = "63-47=16"
$tbl.Cell(10, 2).Range.Text = "8+9=17"
$tbl.Cell(10, 3).Range.Text = "35-7=28"
$tbl.Cell(10, 4).Range.Text = "7+47=54"
$tbl.Cell(10, 5).Range.Text = "39+47=86"
$tbl.Cell(11, 1).Range.Text = "5+88=93"
$tbl.Cell(11, 2).Range.Text = "9+35=44"
$tbl.Cell(11, 3).Range.Text = "80-78=2"
$tbl.Cell(11, 4).Range.Text = "7+64=71"
$tbl.Cell(11, 5).Range.Text = "19+73=92"
$tbl.Cell(12, 1).Range.Text = "47+27=74"
$tbl.Cell(12, 2).Range.Text = "46-27=19"
$tbl.Cell(12, 3).Range.Text = "3+39=42"
$tbl.Cell(12, 4).Range.Text = "82-23=59"
$tbl.Cell(12, 5).Range.Text = "18+77=95"
$tbl.Cell(13, 1).Range.Text = "92-47=45"
$tbl.Cell(13, 2).Range.Text = "25+68=93"
$tbl.Cell(13, 3).Range.Text = "65-8=57"
$tbl.Cell(13, 4).Range.Text = "82-69=13"
$tbl.Cell(13, 5).Range.Text = "44-8=36"
$tbl.Cell(14, 1).Range.Text = "41-14=27"
$tbl.Cell(14, 2).Range.Text = "21-8=13"
$tbl.Cell(14, 3).Range.Text = "61-7=54"
$tbl.Cell(14, 4).Range.Text = "8+57=65"
$tbl.Cell(14, 5).Range.Text = "80-72=8"
$tbl.Cell(15, 1).Range.Text = "19+65=84"
$tbl.Cell(15, 2).Range.Text = "92-15=77"
$tbl.Cell(15, 3).Range.Text = "65+7=72"
$tbl.Cell(15, 4).Range.Text = "49+33=82"
$tbl.Cell(15, 5).Range.Text = "2+49=51"
$tbl.Cell(16, 1).Range.Text = "42-9=33"
$tbl.Cell(16, 2).Range.Text = "94-87=7"
$tbl.Cell(16, 3).Range.Text = "79+15=94"
$tbl.Cell(16, 4).Range.Text = "4+27=31"
$tbl.Cell(16, 5).Range.Text = "69+17=86"
$tbl.Cell(17, 1).Range.Text = "35-29=6"
$tbl.Cell(17, 2).Range.Text = "16+67=83"
$tbl.Cell(17, 3).Range.Text = "86-18=68"
$tbl.Cell(17, 4).Range.Text = "31-22=9"
$tbl.Cell(17, 5).Range.Text = "29+18=47"
$tbl.Cell(18, 1).Range.Text = "56+7=63"
$tbl.Cell(18, 2).Range.Text = "31-16=15"
$tbl.Cell(18, 3).Range.Text = "38-9=29"
$tbl.Cell(18, 4).Range.Text = "33+8=41"
$tbl.Cell(18, 5).Range.Text = "33-19=14"
$tbl.Cell(19, 1).Range.Text = "29+12=41"
$tbl.Cell(19, 2).Range.Text = "52-16=36"
$tbl.Cell(19, 3).Range.Text = "94-38=56"
$tbl.Cell(19, 4).Range.Text = "81-74=7"
$tbl.Cell(19, 5).Range.Text = "37+58=95"
$tbl.Cell(20, 1).Range.Text = "95-46=49"
$tbl.Cell(20, 2).Range.Text = "83-36=47"
$tbl.Cell(20, 3).Range.Text = "74-29=45"
$tbl.Cell(20, 4).Range.Text = "66-7=59"
$tbl.Cell(20, 5).Range.Text = "48+15=63"
